$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-12 Friday", "2024-07-13 Saturday"),
    @("17×65=", "79×67="),
    @("28×85=", "52×66="),
    @("84×55=", "70×44="),
    @("39×74=", "34×50="),
    @("47×56=", "70×36="),
    @("37×62=", "14×25="),
    @("15×92=", "33×79="),
    @("22×92=", "20×47="),
    @("16×58=", "32×57="),
    @("31×84=", "12×90="),
    @("72×21=", "19×24="),
    @("16×45=", "76×24="),
    @("29×97=", "36×80="),
    @("41×37=", "97×88="),
    @("52×18=", "65×25="),
    @("66×58=", "35×91="),
    @("39×47=", "72×53="),
    @("32×95=", "81×40="),
    @("20×69=", "94×96="),
    @("33×18=", "63×21="),
    @("54×56=", "45×33="),
    @("25×30=", "23×77="),
    @("34×47=", "71×90="),
    @("27×26=", "77×22="),
    @("88×92=", "53×99=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
